$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 18875.8
$ws.Range("I51").Value = 5234
$ws.Range("J51").Value = 30812.375
$ws.Range("K51").Value = 5234
$ws.Range("L51").Value = 30812.375
$ws.Range("M51").Value = -4750
$ws.Range("N51").Value = -31780.375

$ws.Range("H52").Value = 853
$ws.Range("I52").Value = 779.5
$ws.Range("J52").Value = 1000
$ws.Range("K52").Value = 2338.5
$ws.Range("L52").Value = 3000
$ws.Range("M52").Value = -2178.5
$ws.Range("N52").Value = -3320

$ws.Range("H64").Value = 6780.9287
$ws.Range("I64").Value = 4704.857
$ws.Range("J64").Value = 8857
$ws.Range("K64").Value = 4704.857
$ws.Range("L64").Value = 8857
$ws.Range("M64").Value = -4456.857
$ws.Range("N64").Value = -9353

$ws.Range("H67").Value = 6780.9287
$ws.Range("I67").Value = 4704.857
$ws.Range("J67").Value = 8857
$ws.Range("K67").Value = 4704.857
$ws.Range("L67").Value = 8857
$ws.Range("M67").Value = -3846.857
$ws.Range("N67").Value = -10573

$ws.Range("H75").Value = 281875
$ws.Range("J75").Value = 281500
$ws.Range("L75").Value = 281500
$ws.Range("N75").Value = -283372

$ws.Range("H78").Value = 281875
$ws.Range("J78").Value = 281500
$ws.Range("L78").Value = 844500
$ws.Range("N78").Value = -853860

$ws.Range("H100").Value = 2311.0588
$ws.Range("J100").Value = 1968.25
$ws.Range("L100").Value = 1968.25
$ws.Range("N100").Value = -3050.25

$ws.Range("H112").Value = 2703.1924
$ws.Range("J112").Value = 2703.1924
$ws.Range("L112").Value = 8109.5772
$ws.Range("N112").Value = -10325.5772

$ws.Range("H132").Value = 1480.4138
$ws.Range("I132").Value = 1169.4
$ws.Range("K132").Value = 3508.2
$ws.Range("M132").Value = -978.2000000000003

$ws.Range("H137").Value = 5561.5
$ws.Range("J137").Value = 8437.036
$ws.Range("L137").Value = 25311.108
$ws.Range("N137").Value = -30411.108

$ws.Range("H138").Value = 3658.724
$ws.Range("J138").Value = 3695.8096
$ws.Range("L138").Value = 11087.4288
$ws.Range("N138").Value = -21367.4288

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 82.05
$ws.Range("I5").Value = 78.47369
$ws.Range("K5").Value = 78.47369
$ws.Range("M5").Value = 33.52631

$ws.Range("H32").Value = 4351192
$ws.Range("I32").Value = 5003257.5
$ws.Range("K32").Value = 5003257.5
$ws.Range("M32").Value = -5002970.5

$ws.Range("H61").Value = 6923.9375
$ws.Range("I61").Value = 4773.75
$ws.Range("K61").Value = 4773.75
$ws.Range("M61").Value = -4561.75

$ws.Range("H110").Value = 21252232
$ws.Range("I110").Value = 34000572
$ws.Range("J110").Value = 5000
$ws.Range("K110").Value = 34000572
$ws.Range("L110").Value = 5000
$ws.Range("M110").Value = -33998527
$ws.Range("N110").Value = -9090

$ws.Range("H132").Value = 8234.541999999999
$ws.Range("I132").Value = 5311.5625
$ws.Range("K132").Value = 15934.6875
$ws.Range("M132").Value = -13404.6875

$ws.Range("H136").Value = 6923.9375
$ws.Range("I136").Value = 4773.75
$ws.Range("K136").Value = 14321.25
$ws.Range("M136").Value = -11771.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 82.05
$ws.Range("I4").Value = 78.47369
$ws.Range("K4").Value = 78.47369
$ws.Range("M4").Value = 36.52631

$ws.Range("H5").Value = 2990.889
$ws.Range("I5").Value = 4238
$ws.Range("K5").Value = 4238
$ws.Range("M5").Value = -4125

$ws.Range("H58").Value = 41690
$ws.Range("I58").Value = 0
$ws.Range("K58").Value = 0
$ws.Range("M58").ClearContents()

$ws.Range("H134").Value = 5826.0527
$ws.Range("I134").Value = 2691.7778
$ws.Range("J134").Value = 8646.9
$ws.Range("K134").Value = 8075.3334
$ws.Range("L134").Value = 25940.7
$ws.Range("M134").Value = -5540.3334
$ws.Range("N134").Value = -31010.7

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2426.0833
$ws.Range("J16").Value = 2850.75
$ws.Range("L16").Value = 2850.75
$ws.Range("N16").Value = -3424.75

$ws.Range("H31").Value = 4067.5405
$ws.Range("I31").Value = 2445.1667
$ws.Range("J31").Value = 5604.5264
$ws.Range("K31").Value = 2445.1667
$ws.Range("L31").Value = 5604.5264
$ws.Range("M31").Value = -2150.1667
$ws.Range("N31").Value = -6194.5264

$ws.Range("H34").Value = 4067.5405
$ws.Range("I34").Value = 2445.1667
$ws.Range("J34").Value = 5604.5264
$ws.Range("K34").Value = 2445.1667
$ws.Range("L34").Value = 5604.5264
$ws.Range("M34").Value = -2243.1667
$ws.Range("N34").Value = -6008.5264

$ws.Range("H41").Value = 4999.1665
$ws.Range("I41").Value = 4999.1665
$ws.Range("J41").Value = 0
$ws.Range("K41").Value = 4999.1665
$ws.Range("L41").Value = 0
$ws.Range("M41").ClearContents()
$ws.Range("N41").Value = -4571.1665

$ws.Range("H113").Value = 2426.0833
$ws.Range("J113").Value = 2850.75
$ws.Range("L113").Value = 2850.75
$ws.Range("N113").Value = -7190.75

$ws.Range("H122").Value = 94133.73
$ws.Range("I122").Value = 146289.14
$ws.Range("K122").Value = 438867.42
$ws.Range("M122").Value = -436417.42

$ws.Range("H132").Value = 53385.16
$ws.Range("I132").Value = 5027.3
$ws.Range("J132").Value = 85623.734
$ws.Range("K132").Value = 15081.9
$ws.Range("L132").Value = 256871.202
$ws.Range("M132").Value = -12551.9
$ws.Range("N132").Value = -261931.202

$ws.Range("H134").Value = 9606.429
$ws.Range("I134").Value = 10140.272
$ws.Range("J134").Value = 7649
$ws.Range("K134").Value = 30420.816
$ws.Range("L134").Value = 22947
$ws.Range("M134").Value = -27885.816
$ws.Range("N134").Value = -28017

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 273.5
$ws.Range("I7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("M7").ClearContents()

$ws.Range("H98").Value = 715.75
$ws.Range("I98").Value = 649.5
$ws.Range("J98").Value = 782
$ws.Range("K98").Value = 1948.5
$ws.Range("L98").Value = 2346
$ws.Range("M98").Value = -450.5
$ws.Range("N98").Value = -5342

$ws.Range("H122").Value = 66671870
$ws.Range("I122").Value = 125009480
$ws.Range("K122").Value = 1125085320
$ws.Range("M122").Value = -1125082870

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 2541.4
$ws.Range("I43").Value = 2541.4
$ws.Range("K43").Value = 2541.4
$ws.Range("M43").Value = -2390.4

$ws.Range("H46").Value = 40860.332
$ws.Range("J46").Value = 49915.5
$ws.Range("L46").Value = 49915.5
$ws.Range("N46").Value = -50227.5

$ws.Range("H97").Value = 2010.6666
$ws.Range("I97").Value = 2031.7
$ws.Range("K97").Value = 2031.7
$ws.Range("M97").Value = -1535.7

$ws.Range("H123").Value = 39174
$ws.Range("J123").Value = 39174
$ws.Range("L123").Value = 39174
$ws.Range("N123").Value = -44074

$ws.Range("H126").Value = 3290.4285
$ws.Range("I126").Value = 2185.3
$ws.Range("J126").Value = 6053.25
$ws.Range("K126").Value = 6555.900000000001
$ws.Range("L126").Value = 18159.75
$ws.Range("M126").Value = -4085.900000000001
$ws.Range("N126").Value = -23099.75

$ws.Range("H132").Value = 8524.861999999999
$ws.Range("I132").Value = 6650.439
$ws.Range("J132").Value = 13045.529
$ws.Range("K132").Value = 19951.317
$ws.Range("L132").Value = 39136.587
$ws.Range("M132").Value = -17421.317
$ws.Range("N132").Value = -44196.587

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("L40").ClearContents()
$ws.Range("M40").ClearContents()
$ws.Range("N40").Value = 0

$ws.Range("H61").Value = 3296.2903
$ws.Range("I61").Value = 3353.2693
$ws.Range("K61").Value = 3353.2693
$ws.Range("M61").Value = -3151.2693

$ws.Range("H113").Value = 3296.2903
$ws.Range("I113").Value = 3353.2693
$ws.Range("K113").Value = 3353.2693
$ws.Range("M113").Value = -1183.2693

$ws.Range("H132").Value = 6472.8667
$ws.Range("I132").Value = 5099.4546
$ws.Range("K132").Value = 15298.3638
$ws.Range("M132").Value = -12768.3638

$ws.Range("H136").Value = 6223.647
$ws.Range("I136").Value = 2756
$ws.Range("K136").Value = 8268
$ws.Range("M136").Value = -5718

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H69").Value = 13445.714
$ws.Range("J69").Value = 14565.6
$ws.Range("L69").Value = 14565.6
$ws.Range("N69").Value = -16063.6

$ws.Range("H72").Value = 13445.714
$ws.Range("J72").Value = 14565.6
$ws.Range("L72").Value = 43696.8
$ws.Range("N72").Value = -51184.8

$ws.Range("H100").Value = 6012.222
$ws.Range("I100").Value = 10510
$ws.Range("K100").Value = 21020
$ws.Range("M100").Value = -20479

$ws.Range("H113").Value = 492.63635
$ws.Range("I113").Value = 504.14285
$ws.Range("J113").Value = 472.5
$ws.Range("K113").Value = 1512.42855
$ws.Range("L113").Value = 1417.5
$ws.Range("M113").Value = 657.5714499999999
$ws.Range("N113").Value = -5757.5

$ws.Range("H126").Value = 12628407
$ws.Range("I126").Value = 16836218
$ws.Range("K126").Value = 50508654
$ws.Range("M126").Value = -50506184

$ws.Range("H132").Value = 4024.1914
$ws.Range("I132").Value = 1824.9354
$ws.Range("K132").Value = 5474.8062
$ws.Range("M132").Value = -2944.8062

$ws.Range("H136").Value = 4329.243
$ws.Range("I136").Value = 3252.7856
$ws.Range("J136").Value = 7678.222
$ws.Range("K136").Value = 9758.356800000001
$ws.Range("L136").Value = 23034.666
$ws.Range("M136").Value = -7208.356800000001
$ws.Range("N136").Value = -28134.666
